$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Add()
$vals = @(0.001234567890123456, 0.01234567890123456, 0.1234567890123456, 1.234567890123456, 12.34567890123456, 123.4567890123456, 1234.567890123456, 12345.67890123456)
for ($i=0; $i -lt $vals.Length; $i++) {
  $r = 1 + $i
  $ws2.Cells.Item($r, 1).Value2 = $vals[$i]
}
